# Fruta / hortaliza, semanal
# Insert a new weekly price-report row at row 78 (pushing all existing
# rows 78..184 down to 79..185) and populate it with the new week's data.
# All other rows keep their original contents; they merely shift down one
# position as a side effect of the insert.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 78; formatting (e.g. the date
# number format on column D) is inherited from the row above, matching
# how the rest of the sheet is styled.
$ws.Rows.Item(78).Insert()

$ws.Cells.Item(78, 1).Value = 3
$ws.Cells.Item(78, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(78, 3).Value = "Coquimbo"
$ws.Cells.Item(78, 4).Value = 44413
$ws.Cells.Item(78, 5).Value = 5
$ws.Cells.Item(78, 6).Value = 100112003
$ws.Cells.Item(78, 7).Value = "Ajo"
$ws.Cells.Item(78, 8).Value = "Chino"
$ws.Cells.Item(78, 9).Value = "Primera"
$ws.Cells.Item(78, 10).Value = 95
$ws.Cells.Item(78, 11).Value = 11000
$ws.Cells.Item(78, 12).Value = 11500
$ws.Cells.Item(78, 13).Value = 11237
$ws.Cells.Item(78, 14).Value = "`$/caja 10 kilos"
$ws.Cells.Item(78, 15).Value = "China"
$ws.Cells.Item(78, 16).Value = 1124
$ws.Cells.Item(78, 17).Value = 10
$ws.Cells.Item(78, 18).Value = "Hortaliza"
